$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 7500
$ws.Range("I21").Value = 5000
$ws.Range("J21").Value = 10000
$ws.Range("K21").Value = 5000
$ws.Range("L21").Value = 10000
$ws.Range("M21").Value = -4532
$ws.Range("N21").Value = -10936

$ws.Range("H23").Value = 7500
$ws.Range("I23").Value = 5000
$ws.Range("J23").Value = 10000
$ws.Range("K23").Value = 5000
$ws.Range("L23").Value = 10000
$ws.Range("M23").Value = -4766
$ws.Range("N23").Value = -10468

$ws.Range("H74").Value = 13794.25
$ws.Range("I74").Value = 11725.667
$ws.Range("K74").Value = 11725.667
$ws.Range("M74").Value = -10789.667

$ws.Range("H77").Value = 13794.25
$ws.Range("I77").Value = 11725.667
$ws.Range("K77").Value = 58628.335
$ws.Range("M77").Value = -53948.335

$ws.Range("H98").Value = 1965
$ws.Range("I98").Value = 1948
$ws.Range("K98").Value = 1948
$ws.Range("M98").Value = -450

$ws.Range("H122").Value = 1965
$ws.Range("I122").Value = 1948
$ws.Range("K122").Value = 5844
$ws.Range("M122").Value = -3394

$ws.Range("H125").Value = 1526.2
$ws.Range("I125").Value = 1444
$ws.Range("J125").Value = 1649.5
$ws.Range("K125").Value = 12996
$ws.Range("L125").Value = 14845.5
$ws.Range("M125").Value = -10536
$ws.Range("N125").Value = -19765.5

$ws.Range("H131").Value = 945
$ws.Range("I131").Value = 945
$ws.Range("K131").Value = 2835
$ws.Range("M131").Value = 2205

$ws.Range("H137").Value = 2805.5518
$ws.Range("I137").Value = 2018.3
$ws.Range("J137").Value = 4555
$ws.Range("K137").Value = 6054.9
$ws.Range("L137").Value = 13665
$ws.Range("M137").Value = -3504.9
$ws.Range("N137").Value = -18765

$ws.Range("H141").Value = 2748.75
$ws.Range("I141").Value = 965
$ws.Range("K141").Value = 2895
$ws.Range("M141").Value = 2285

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1899.4
$ws.Range("I74").Value = 1561.3793
$ws.Range("J74").Value = 3533.1667
$ws.Range("K74").Value = 1561.3793
$ws.Range("L74").Value = 3533.1667
$ws.Range("M74").Value = -687.3793000000001
$ws.Range("N74").Value = -5281.1667

$ws.Range("H77").Value = 1899.4
$ws.Range("I77").Value = 1561.3793
$ws.Range("J77").Value = 3533.1667
$ws.Range("K77").Value = 7806.896500000001
$ws.Range("L77").Value = 17665.8335
$ws.Range("M77").Value = -3438.896500000001
$ws.Range("N77").Value = -26401.8335

$ws.Range("H110").Value = 704.5
$ws.Range("I110").Value = 445.6
$ws.Range("K110").Value = 445.6
$ws.Range("M110").Value = 1599.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H97").Value = 16941
$ws.Range("I97").Value = 16941
$ws.Range("K97").Value = 16941
$ws.Range("M97").Value = -15950

$ws.Range("H105").Value = 2185
$ws.Range("I105").Value = 2239.818
$ws.Range("J105").Value = 2064.4
$ws.Range("K105").Value = 2239.818
$ws.Range("L105").Value = 2064.4
$ws.Range("M105").Value = -492.8180000000002
$ws.Range("N105").Value = -5558.4

$ws.Range("H107").Value = 3854.6875
$ws.Range("I107").Value = 3093.625
$ws.Range("K107").Value = 3093.625
$ws.Range("M107").Value = -1173.625

$ws.Range("H134").Value = 2419.3
$ws.Range("I134").Value = 1872.4667
$ws.Range("K134").Value = 5617.4001
$ws.Range("M134").Value = -3082.4001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 390.625
$ws.Range("I19").Value = 160.71428
$ws.Range("K19").Value = 160.71428
$ws.Range("M19").Value = 9.285719999999998

$ws.Range("H24").Value = 390.625
$ws.Range("I24").Value = 160.71428
$ws.Range("K24").Value = 160.71428
$ws.Range("M24").Value = 9.285719999999998

$ws.Range("H132").Value = 2414.8
$ws.Range("I132").Value = 2628.2222
$ws.Range("K132").Value = 7884.6666
$ws.Range("M132").Value = -5354.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 984
$ws.Range("I3").Value = 984
$ws.Range("K3").Value = 2952
$ws.Range("M3").Value = -2840

$ws.Range("H68").Value = 1546.2354
$ws.Range("I68").Value = 1365.4
$ws.Range("J68").Value = 1621.5834
$ws.Range("K68").Value = 4096.200000000001
$ws.Range("L68").Value = 4864.7502
$ws.Range("M68").Value = -3285.200000000001
$ws.Range("N68").Value = -6486.7502

$ws.Range("H71").Value = 1546.2354
$ws.Range("I71").Value = 1365.4
$ws.Range("J71").Value = 1621.5834
$ws.Range("K71").Value = 12288.6
$ws.Range("L71").Value = 14594.2506
$ws.Range("M71").Value = -8232.6
$ws.Range("N71").Value = -22706.2506

$ws.Range("H103").Value = 563
$ws.Range("I103").Value = 268.75
$ws.Range("J103").Value = 759.1667
$ws.Range("K103").Value = 806.25
$ws.Range("L103").Value = 2277.5001
$ws.Range("M103").Value = 72.75
$ws.Range("N103").Value = -4035.5001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 29999
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 29999
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 29999
$ws.Range("M47").ClearContents()
$ws.Range("N47").Value = -31135

$ws.Range("H97").Value = 987
$ws.Range("I97").Value = 925
$ws.Range("K97").Value = 925
$ws.Range("M97").Value = -429

$ws.Range("H126").Value = 2082.3333
$ws.Range("I126").Value = 2082.3333
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 6246.999899999999
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -3776.999899999999
$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 3171
$ws.Range("I132").Value = 2590
$ws.Range("J132").Value = 4236.1665
$ws.Range("K132").Value = 7770
$ws.Range("L132").Value = 12708.4995
$ws.Range("M132").Value = -5240
$ws.Range("N132").Value = -17768.4995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H92").Value = 26666.666
$ws.Range("I92").Value = 10000
$ws.Range("K92").Value = 10000
$ws.Range("M92").Value = -7504

$ws.Range("H93").Value = 1000
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 5397.85
$ws.Range("I126").Value = 2907.889
$ws.Range("J126").Value = 7435.091
$ws.Range("K126").Value = 8723.667000000001
$ws.Range("L126").Value = 22305.273
$ws.Range("M126").Value = -6253.667000000001
$ws.Range("N126").Value = -27245.273

$ws.Range("H132").Value = 2053.1538
$ws.Range("I132").Value = 1710.6666
$ws.Range("K132").Value = 5131.9998
$ws.Range("M132").Value = -2601.9998

$ws.Range("H136").Value = 3277
$ws.Range("I136").Value = 2495.0588
$ws.Range("J136").Value = 4606.3
$ws.Range("K136").Value = 7485.176399999999
$ws.Range("L136").Value = 13818.9
$ws.Range("M136").Value = -4935.176399999999
$ws.Range("N136").Value = -18918.9
